# The commit scrolled/re-selected Sheet1 and saved the workbook, which (because
# column B is full of volatile CHAR(RANDBETWEEN(...)) formulas marked ca="1")
# re-rolled every cached formula result in the process. Reproduce both the
# view-state change and the recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Scroll the window so row 58 is pinned at the top (topLeftCell="A58" in the
# saved view) before moving the selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1

# New selection: activeCell="A79" sqref="A79"
$ws.Range("A79").Select() | Out-Null

# Force a full recalculation so every "=CHAR(RANDBETWEEN(65,90))&..." cell in
# column B (B2, B4, B21, B23, ... B187) gets a freshly generated cached value,
# matching the bulk of <v> changes in the diff.
$excel.CalculateFull()
